# Apply cryptos list update (cell value changes) via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.180.48"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "2.359.59"
$ws.Range("E3").Value = "  +0.93%  "
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.42%  "
$c = $ws.Range("D5")
$c.Value = "'542.01"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.57%  "
$c = $ws.Range("D6")
$c.Value = "'136.45"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.09%  "
$c = $ws.Range("D7")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "
$c = $ws.Range("D8")
$c.Value = "'0.564"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +5.46%  "
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("E10").Value = "  +4.59%  "
$ws.Range("E11").Value = "  -0.65%  "
$c = $ws.Range("D12")
$c.Value = "'0.355"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.98%  "
$c = $ws.Range("D13")
$c.Value = "'23.95"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "2.774.96"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "58.126.39"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "2.344.15"
$ws.Range("E17").Value = "  +1.06%  "
$c = $ws.Range("D18")
$c.Value = "'10.74"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.13%  "
$c = $ws.Range("D19")
$c.Value = "'332.78"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "
$c = $ws.Range("D20")
$c.Value = "'4.29"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.76%  "
$c = $ws.Range("D21")
$c.Value = "'6.75"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("E22").Value = "  +0.24%  "
$c = $ws.Range("D23")
$c.Value = "'62.93"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D25")
$c.Value = "'8.54"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D26")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "
$c = $ws.Range("D27")
$c.Value = "'1.39"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.72%  "
$c = $ws.Range("D28")
$c.Value = "'172.63"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("D30").Value = "0.0₃0741"
$ws.Range("E30").Value = "  +2.46%  "
$c = $ws.Range("D31")
$c.Value = "'6.17"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  +12.09%  "
$c = $ws.Range("D33")
$c.Value = "'18.53"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D35")
$c.Value = "'4.24"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +6.78%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D36")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("E38").Value = "  +4.52%  "
$c = $ws.Range("D39")
$c.Value = "'39.34"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.41%  "
$c = $ws.Range("D40")
$c.Value = "'145.59"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "
$c = $ws.Range("D41")
$c.Value = "'293.84"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("E42").Value = "  +1.04%  "
$c = $ws.Range("D43")
$c.Value = "'3.66"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "
$c = $ws.Range("D44")
$c.Value = "'0.0948"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.72%  "
$c = $ws.Range("D45")
$c.Value = "'19.38"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.39%  "
$c = $ws.Range("D46")
$c.Value = "'0.0504"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "
$c = $ws.Range("D47")
$c.Value = "'0.566"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.56%  "
$c = $ws.Range("D48")
$c.Value = "'0.0222"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D49")
$c.Value = "'17.51"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D50")
$c.Value = "'0.382"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$c = $ws.Range("D51")
$c.Value = "'11.07"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
